# "new features to bokeh server" edit
# - investment toggle button added -> drop the old "potential restriction (MW_th)"
#   column (J) on the Parameter sheet, shifting the "renewable factor" column
#   (formerly K) left into J, and refresh the installed-capacity (C) figures.
# - "heat profile" selection / "electricity-price" selection -> shared-string
#   bookkeeping is handled automatically by the runtime when the column is
#   removed, so sheet2/sheet3 need no direct edits.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Remove the "potential restriction (MW_th)" column entirely; this shifts the
# old column K ("renewable factor") into column J and compacts the shared
# string table automatically.
$ws1.Columns("J").Delete()

# Refresh the "installed capacity (MW_th)" values (column C) for every power
# plant row with the new figures.
$capacities = @(50, 100, 150, 200, 250, 300, 350, 400, 450, 500, 550, 600, 650, 750, 850, 900, 1000, 1050, 1100, 1150)
for ($i = 0; $i -lt $capacities.Length; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 3).Value = $capacities[$i]
}

# Row 3 ("renewable factor", now column J) changes from 0.1 to 1.
$ws1.Range("J3").Value = 1

# Minor floating point refresh on the OPEX var figure for row 11.
$ws1.Range("H11").Value = 0.8999999999999999
